$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'51.276.97"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -1.66%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.919.97"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -2.45%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D5").Value = "'372.74"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +5.16%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'103.28"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -4.16%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  -3.74%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  -0.08%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.590"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -5.95%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'37.10"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -3.14%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  +0.42%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  -2.29%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'  -5.19%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'3.379.52"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -2.39%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'7.40"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -4.03%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'2.914.78"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -2.77%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Value = "'  -8.69%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'51.224.07"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -1.88%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'3.31"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -5.15%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'7.28"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -3.44%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'12.99"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -5.07%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.0₃0948"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -2.57%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'68.44"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -1.57%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'260.83"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -1.23%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'2.69"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -1.42%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.171"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -4.64%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("B27").Value = "'LEO"
$ws.Range("B27").Style = "Normal"
$ws.Range("C27").Value = "'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("C27").Style = "Normal"
$ws.Range("D27").Value = "'4.12"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -4.16%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("B28").Value = "'Dai"
$ws.Range("B28").Style = "Normal"
$ws.Range("C28").Value = "'https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("C28").Style = "Normal"
$ws.Range("D28").Value = "'1.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -0.01%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'25.79"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -4.24%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'7.17"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -5.89%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'6.52"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +5.32%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  -6.40%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'9.90"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -3.98%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  -3.65%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'34.66"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -5.02%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'51.23"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +0.72%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  +0.46%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  -4.22%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  -6.56%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'17.11"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -4.37%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'2.60"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -4.20%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'1.85"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -6.39%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'  -3.54%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'22.11"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -3.12%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'119.78"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -2.09%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'  -2.83%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'2.024.49"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -4.72%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'  -3.90%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'  -6.09%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.245"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +1.60%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'3.208.74"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -2.34%  "
$ws.Range("E51").Style = "Normal"
